$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, pushing existing rows 9-27 down to 10-28
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the weekly record (same constants as the rest of
# this dataset; only Fecha/Volumen/Precio columns differ row to row).
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 45070
$ws.Cells.Item(9, 4).NumberFormat = $ws.Cells.Item(10, 4).NumberFormat
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = 100112039
$ws.Cells.Item(9, 7).Value = "Ciboulette"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 800
$ws.Cells.Item(9, 11).Value = 2000
$ws.Cells.Item(9, 12).Value = 2500
$ws.Cells.Item(9, 13).Value = 2250
$ws.Cells.Item(9, 14).Value = '$/docena de atados'
$ws.Cells.Item(9, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(9, 16).Value = 750
$ws.Cells.Item(9, 17).Value = 3
$ws.Cells.Item(9, 18).Value = "Hortaliza"
